$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B1").Value = 16
$ws.Range("C1").Value = 20
$ws.Range("D1").Value = 16
$ws.Range("E1").Value = 20

$ws.Range("B2").Value = 46.75668948809188
$ws.Range("C2").Value = 53.719528762676511
$ws.Range("D2").Value = 49.966250868235662
$ws.Range("E2").Value = 54.39359641041851

$ws.Range("B3").Value = 44.067825498757585
$ws.Range("C3").Value = 44.125118230966279
$ws.Range("D3").Value = 43.591477932075151
$ws.Range("E3").Value = 53.34815552156919

$ws.Range("B1:E3").Select()
